$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set C8 to empty string (value present, still inlineStr with no text)
$ws.Cells.Item(8, 3).Value = ""

# Add new row 9 with the new annotation data
$ws.Cells.Item(9, 1).Value = "parisk"
$ws.Cells.Item(9, 2).Value = 3
$ws.Cells.Item(9, 3).Value = "nan"
$ws.Cells.Item(9, 4).Value = "SUG"
$ws.Cells.Item(9, 5).Value = "WRI"
$ws.Cells.Item(9, 6).Value = "41c93df3-3a59-4ce4-b94b-f420b7540586"
$ws.Cells.Item(9, 7).Value = "SJ19eUg0-_annotated.xlsx"
$ws.Cells.Item(9, 8).Value = "Thank the reviewer for the thoughtful feedback."
